$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" date column (C) for rows 2-8 from serial 45174 to 45175
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 3).Value = 45175
}
